$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6583
$ws.Range("J3").Value = 6967
$ws.Range("C4").Value = 1838
$ws.Range("J4").Value = 1513
$ws.Range("J5").Value = 540
$ws.Range("J6").Value = 9263
$ws.Range("C7").Value = 28382
$ws.Range("J7").Value = 24866

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 429
$ws.Range("J3").Value = 469
$ws.Range("J6").Value = 548
$ws.Range("J7").Value = 1567

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 372
$ws.Range("J4").Value = 51
$ws.Range("J6").Value = 394
$ws.Range("J7").Value = 1125

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 223
$ws.Range("J3").Value = 258
$ws.Range("J4").Value = 27
$ws.Range("J6").Value = 224
$ws.Range("J7").Value = 762

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 182
$ws.Range("J3").Value = 174
$ws.Range("J6").Value = 224
$ws.Range("J7").Value = 623

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 105
$ws.Range("J3").Value = 150
$ws.Range("J7").Value = 378

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J8").Value = 1567
$ws.Range("J9").Value = 133
$ws.Range("J10").Value = 182
$ws.Range("J11").Value = 426
$ws.Range("J14").Value = 134
$ws.Range("J15").Value = 293
$ws.Range("J19").Value = 727
$ws.Range("J20").Value = 521
$ws.Range("J22").Value = 59
$ws.Range("J27").Value = 149
$ws.Range("J29").Value = 1356
$ws.Range("J31").Value = 244
$ws.Range("J33").Value = 1125
$ws.Range("J37").Value = 762
$ws.Range("J41").Value = 173
$ws.Range("J42").Value = 1071
$ws.Range("J44").Value = 188
$ws.Range("J47").Value = 185
$ws.Range("J48").Value = 281
$ws.Range("J52").Value = 630
$ws.Range("J54").Value = 468
$ws.Range("J55").Value = 382
$ws.Range("J57").Value = 112
$ws.Range("C63").Value = 268
$ws.Range("J63").Value = 82
$ws.Range("J65").Value = 623
$ws.Range("J67").Value = 941
$ws.Range("J69").Value = 54
$ws.Range("J71").Value = 82
$ws.Range("J76").Value = 368
$ws.Range("J79").Value = 696
$ws.Range("J85").Value = 1032
$ws.Range("J90").Value = 264
$ws.Range("J91").Value = 285
$ws.Range("J94").Value = 261
$ws.Range("J96").Value = 274
$ws.Range("J97").Value = 224
$ws.Range("J99").Value = 378
$ws.Range("C101").Value = 28382
$ws.Range("J101").Value = 24866

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J4").Value = 14
$ws.Range("J5").Value = 5
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 240
$ws.Range("J3").Value = 350
$ws.Range("J7").Value = 941

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 98
$ws.Range("J7").Value = 468

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 415
$ws.Range("J3").Value = 475
$ws.Range("J6").Value = 343
$ws.Range("J7").Value = 1356

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 52
$ws.Range("J7").Value = 281

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 176
$ws.Range("J6").Value = 281
$ws.Range("J7").Value = 727

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 62
$ws.Range("J3").Value = 77
$ws.Range("J7").Value = 368

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 173

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 210
$ws.Range("J6").Value = 572
$ws.Range("J7").Value = 1071

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 382

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J2").Value = 83
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 274

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 119
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 285

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 236
$ws.Range("J7").Value = 696

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 147
$ws.Range("J6").Value = 147
$ws.Range("J7").Value = 521

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 51
$ws.Range("J7").Value = 261

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 185

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J6").Value = 126
$ws.Range("J7").Value = 293

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 125
$ws.Range("J6").Value = 191
$ws.Range("J7").Value = 426

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 39
$ws.Range("J6").Value = 155
$ws.Range("J7").Value = 224

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J5").Value = 10
$ws.Range("J7").Value = 264

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 68
$ws.Range("J3").Value = 79

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 277
$ws.Range("J3").Value = 364
$ws.Range("J6").Value = 299
$ws.Range("J7").Value = 1032

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 144
$ws.Range("J6").Value = 271
$ws.Range("J7").Value = 630
